$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in row 36 column A ("pltos" -> "plots") ---
$ws.Range("A36").Value = 'faire une boucle dans "appareiller les plots"'

# --- Fill in "Qui fait quoi ?" (assignee) column B for the action-item rows ---
$ws.Range("B34").Value = "Jordan / CeauMar"
$ws.Range("B35").Value = "Jordan / CeauMar"
$ws.Range("B36").Value = "kiki / charly"
$ws.Range("B37").Value = "kiki / charly"
$ws.Range("B38").Value = "kiki / charly"
$ws.Range("B39").Value = "NZ / KIKI"
$ws.Range("B40").Value = "Jordan / CeauMar"
$ws.Range("B41").Value = "Tout le monde"

# B34 keeps the word-wrap formatting applied on it
$ws.Range("B34").WrapText = $true

# --- Mark some rows as "done" with the green "Satisfaisant" cell style in column C ---
$ws.Range("C34").Style = "Satisfaisant"
$ws.Range("C35").Style = "Satisfaisant"
$ws.Range("C36").Style = "Satisfaisant"
$ws.Range("C37").Style = "Satisfaisant"
$ws.Range("C40").Style = "Satisfaisant"

# --- New action-item row 42 ---
$ws.Range("A42").Value = "Reprendre la maquette état des plots "
$ws.Range("B42").Value = "kiki / charly"
$ws.Range("C42").Style = "Satisfaisant"

# --- Update the view state (scroll position + active cell selection) ---
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("B37").Select()
